# Generate report for IND
# Update forecast values on row 14 (IND) and set active selection on sheet "data"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update row 14 forecast values (columns H through Q)
$ws.Range("H14").Value = -8
$ws.Range("I14").Value = 10.3
$ws.Range("J14").Value = -10
$ws.Range("K14").Value = 11.4
$ws.Range("L14").Value = -15.7
$ws.Range("M14").Value = 15.2
$ws.Range("N14").Value = -9.9
$ws.Range("O14").Value = 12.5
$ws.Range("P14").Value = 6.3
$ws.Range("Q14").Value = 4.7

# Update the active cell selection on the frozen-pane sheet view
$ws.Activate()
$ws.Range("H15").Select()
